# Update the "想去人数" (want-to-go count) figures that changed between
# the previous gh-pages data snapshot and the new one generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 21
$ws1.Range("F7").Value  = 14885
$ws1.Range("F8").Value  = 404
$ws1.Range("F9").Value  = 2
$ws1.Range("F11").Value = 15209
$ws1.Range("F13").Value = 8716
$ws1.Range("F16").Value = 64
$ws1.Range("F22").Value = 508
$ws1.Range("F36").Value = 268
$ws1.Range("F39").Value = 5353

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 60

# --- Sheet "全部类型" (All types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 21
$ws4.Range("F7").Value  = 14886
$ws4.Range("F8").Value  = 404
$ws4.Range("F9").Value  = 2
$ws4.Range("F11").Value = 15209
$ws4.Range("F13").Value = 8716
$ws4.Range("F17").Value = 64
$ws4.Range("F23").Value = 508
$ws4.Range("F33").Value = 60
$ws4.Range("F39").Value = 268
$ws4.Range("F42").Value = 5353
